$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix author name typo in the "Los Juegos del Hambre" row: "Luisa Pedroza" -> "Luisa Pedraza"
$ws.Range("C3").Value = "Luisa Pedraza"

# Scroll the sheet view back to the top (clears the topLeftCell="A3" scroll offset)
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
